$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# Price cells are kept as text (matching source formatting like "27.844.29")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.844.29"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.627.03"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.43"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.28"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.81"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.622.44"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.561"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.30"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.831.29"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.42"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.62"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.07"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.03"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.09"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.393.71"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +10.44%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.846"
$ws.Range("E40").Value = "  -3.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.83"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.57"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.43"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.762.71"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.15"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.84"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.58"
$ws.Range("E51").Value = "  +0.29%  "

Write-Output "Applied cryptos update"
